{"js": "// Update the WMO Met Data Exchange Interoperability Experiment document:\n//   - bump the document date from 2025-01-28 to 2025-02-03 (title-page date\n//     paragraph and the \"Date:\" row in the cover table)\n//   - sentence-case three Heading1 titles (\"Executive Summary\",\n//     \"High Level Architecture\", \"Revision History\")\n\nconst body = context.document.body;\n\nasync function replaceAll(searchText, replacement, matchCase) {\n  const results = body.search(searchText, { matchCase: matchCase, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Exact, case-sensitive replacements so we only touch the intended text.\nawait replaceAll(\"Date: 2025-01-28\", \"Date: 2025-02-03\", true);\nawait replaceAll(\"2025-01-28\", \"2025-02-03\", true);\nawait replaceAll(\"Executive Summary\", \"Executive summary\", true);\nawait replaceAll(\"High Level Architecture\", \"High level architecture\", true);\nawait replaceAll(\"Revision History\", \"Revision history\", true);\n", "ps1": "# Update the WMO Met Data Exchange Interoperability Experiment document:\n#   - bump the document date from 2025-01-28 to 2025-02-03 (covers both the\n#     title-page date paragraph and the \"Date:\" row in the cover table)\n#   - sentence-case three Heading1 titles (\"Executive Summary\",\n#     \"High Level Architecture\", \"Revision History\")\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-AllText \"2025-01-28\" \"2025-02-03\"\nReplace-AllText \"Executive Summary\" \"Executive summary\"\nReplace-AllText \"High Level Architecture\" \"High level architecture\"\nReplace-AllText \"Revision History\" \"Revision history\"\n"}
